# Added size of parts and station
# Adds two small reference tables (Part radius/height, and P Station height)
# to the right of the existing dimensions table, with a bold/blue header
# style and thin borders on the data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values (order matters: it controls the order new shared
#    strings get interned in xl/sharedStrings.xml)
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Part"
$ws.Range("F2").Value = "Radius"
$ws.Range("G1").Value = "(mm)"
$ws.Range("F3").Value = "Height"
$ws.Range("F5").Value = "P Station"
$ws.Range("G5").Value = "(mm)"
$ws.Range("F6").Value = "Height"

$ws.Range("G2").Value = 20
$ws.Range("G3").Value = 25
$ws.Range("G6").Value = 900

# ---------------------------------------------------------------------
# 2. Formatting
#    Build each unique style once on a single cell, then replicate it
#    with Copy/PasteSpecial(formats) so the same style index is reused
#    everywhere instead of minting a fresh one per cell.
# ---------------------------------------------------------------------

# Style A: header row - bold white text on a blue (accent1) fill with a
# thin border all around.
$hdr = $ws.Range("F1")
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2
$hdr.Interior.ThemeColor = 5
$hdr.Font.Bold = $true
$hdr.Font.ThemeColor = 2
$hdr.HorizontalAlignment = 1
$hdr.Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)

# Style B: label cell ("Radius") - plain Calibri font with a thin border.
$lbl = $ws.Range("F2")
$lbl.Borders.LineStyle = 1
$lbl.Borders.Weight = 2
$lbl.Font.Name = "Calibri"

# Style C: plain value/label cells with just a thin border.
$plain = $ws.Range("G2")
$plain.Borders.LineStyle = 1
$plain.Borders.Weight = 2
$plain.Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("G6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Page setup (paper size / orientation)
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. Restore the on-screen selection to match the saved workbook
# ---------------------------------------------------------------------
$ws.Range("I10").Select()
